$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 and Row 19 swap content (Polkadot <-> TRON) plus updated price/volume
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.60"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.121"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.07%  "

# Remaining rows: update Price (D) and/or Volume(1h) (E) values

$ws.Range("D2").Value = "70.191.93"
$ws.Range("E2").Value = "  -0.08%  "

$ws.Range("D3").Value = "3.747.87"
$ws.Range("E3").Value = "  -1.40%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "620.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.98%  "

$ws.Range("D7").Value = "3.746.33"
$ws.Range("E7").Value = "  -1.49%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -3.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.38"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.483"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.28"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.19%  "

$ws.Range("E14").Value = "  -1.31%  "

$ws.Range("D15").Value = "4.370.24"
$ws.Range("E15").Value = "  -1.43%  "

$ws.Range("D16").Value = "3.750.10"
$ws.Range("E16").Value = "  -1.24%  "

$ws.Range("D17").Value = "70.224.60"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "505.37"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.38%  "

$ws.Range("E22").Value = "  -3.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.723"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.55"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.11"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.24%  "

$ws.Range("E27").Value = "  +1.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000133"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.72%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.48"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.06%  "

$ws.Range("E31").Value = "  +1.56%  "

$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.53"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.40%  "

$ws.Range("E34").Value = "  -0.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.12"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.349"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.140"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.22"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +16.29%  "

$ws.Range("E41").Value = "  -4.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.91"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "431.66"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.49"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.63"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.98%  "

$ws.Range("D46").Value = "2.965.20"
$ws.Range("E46").Value = "  -5.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0363"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.28"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.24"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.20%  "

$ws.Range("E51").Value = "  -2.96%  "
